$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PSSM score matrix (B2:K21) with recomputed values
# from the supplemental-figure re-run.
$ws.Range("B2").Value = -19.23066314620145
$ws.Range("C2").Value = 1.979675853413211
$ws.Range("D2").Value = -19.23066314620145
$ws.Range("E2").Value = -19.23066314620145
$ws.Range("F2").Value = -19.23066314620145
$ws.Range("G2").Value = -19.23066314620145
$ws.Range("H2").Value = -19.23066314620145
$ws.Range("I2").Value = -19.23066314620145
$ws.Range("J2").Value = -19.23066314620145
$ws.Range("K2").Value = -19.23066314620145
$ws.Range("B3").Value = -19.23066314620145
$ws.Range("C3").Value = -19.23066314620145
$ws.Range("D3").Value = -19.23066314620145
$ws.Range("E3").Value = -19.23066314620145
$ws.Range("F3").Value = -19.23066314620145
$ws.Range("G3").Value = -19.23066314620145
$ws.Range("H3").Value = -19.23066314620145
$ws.Range("I3").Value = -19.23066314620145
$ws.Range("J3").Value = -19.23066314620145
$ws.Range("K3").Value = -19.23066314620145
$ws.Range("B4").Value = -19.23066314620145
$ws.Range("C4").Value = 1.983381406256167
$ws.Range("D4").Value = 1.684204429072244
$ws.Range("E4").Value = -19.23066314620145
$ws.Range("F4").Value = 3.406206813090976
$ws.Range("G4").Value = -19.23066314620145
$ws.Range("H4").Value = 1.613188149184301
$ws.Range("I4").Value = -19.23066314620145
$ws.Range("J4").Value = 1.105417308296639
$ws.Range("K4").Value = -19.23066314620145
$ws.Range("B5").Value = -19.23066314620145
$ws.Range("C5").Value = 1.661091637572716
$ws.Range("D5").Value = -19.23066314620145
$ws.Range("E5").Value = -19.23066314620145
$ws.Range("F5").Value = -19.23066314620145
$ws.Range("G5").Value = 2.847254675436115
$ws.Range("H5").Value = -19.23066314620145
$ws.Range("I5").Value = -19.23066314620145
$ws.Range("J5").Value = -19.23066314620145
$ws.Range("K5").Value = -19.23066314620145
$ws.Range("B6").Value = -19.23066314620145
$ws.Range("C6").Value = -19.23066314620145
$ws.Range("D6").Value = -19.23066314620145
$ws.Range("E6").Value = -19.23066314620145
$ws.Range("F6").Value = -19.23066314620145
$ws.Range("G6").Value = -19.23066314620145
$ws.Range("H6").Value = -19.23066314620145
$ws.Range("I6").Value = -19.23066314620145
$ws.Range("J6").Value = -19.23066314620145
$ws.Range("K6").Value = -19.23066314620145
$ws.Range("B7").Value = 2.455938683555104
$ws.Range("C7").Value = -19.23066314620145
$ws.Range("D7").Value = -19.23066314620145
$ws.Range("E7").Value = -19.23066314620145
$ws.Range("F7").Value = -19.23066314620145
$ws.Range("G7").Value = -19.23066314620145
$ws.Range("H7").Value = -19.23066314620145
$ws.Range("I7").Value = -19.23066314620145
$ws.Range("J7").Value = -19.23066314620145
$ws.Range("K7").Value = -19.23066314620145
$ws.Range("B8").Value = -19.23066314620145
$ws.Range("C8").Value = -19.23066314620145
$ws.Range("D8").Value = -19.23066314620145
$ws.Range("E8").Value = 1.795763711058337
$ws.Range("F8").Value = -19.23066314620145
$ws.Range("G8").Value = -19.23066314620145
$ws.Range("H8").Value = -19.23066314620145
$ws.Range("I8").Value = -19.23066314620145
$ws.Range("J8").Value = -19.23066314620145
$ws.Range("K8").Value = -19.23066314620145
$ws.Range("B9").Value = 3.859300783889804
$ws.Range("C9").Value = -19.23066314620145
$ws.Range("D9").Value = -19.23066314620145
$ws.Range("E9").Value = -19.23066314620145
$ws.Range("F9").Value = -19.23066314620145
$ws.Range("G9").Value = -19.23066314620145
$ws.Range("H9").Value = -19.23066314620145
$ws.Range("I9").Value = -19.23066314620145
$ws.Range("J9").Value = -19.23066314620145
$ws.Range("K9").Value = -19.23066314620145
$ws.Range("B10").Value = -19.23066314620145
$ws.Range("C10").Value = -19.23066314620145
$ws.Range("D10").Value = -19.23066314620145
$ws.Range("E10").Value = -19.23066314620145
$ws.Range("F10").Value = -19.23066314620145
$ws.Range("G10").Value = -19.23066314620145
$ws.Range("H10").Value = -19.23066314620145
$ws.Range("I10").Value = -19.23066314620145
$ws.Range("J10").Value = -19.23066314620145
$ws.Range("K10").Value = 2.191116632116316
$ws.Range("B11").Value = -19.23066314620145
$ws.Range("C11").Value = -19.23066314620145
$ws.Range("D11").Value = -19.23066314620145
$ws.Range("E11").Value = 2.90866742700753
$ws.Range("F11").Value = -19.23066314620145
$ws.Range("G11").Value = 2.810507024527246
$ws.Range("H11").Value = -19.23066314620145
$ws.Range("I11").Value = -19.23066314620145
$ws.Range("J11").Value = -19.23066314620145
$ws.Range("K11").Value = 1.879958432595656
$ws.Range("B12").Value = -19.23066314620145
$ws.Range("C12").Value = -19.23066314620145
$ws.Range("D12").Value = -19.23066314620145
$ws.Range("E12").Value = -19.23066314620145
$ws.Range("F12").Value = -19.23066314620145
$ws.Range("G12").Value = -19.23066314620145
$ws.Range("H12").Value = -19.23066314620145
$ws.Range("I12").Value = -19.23066314620145
$ws.Range("J12").Value = -19.23066314620145
$ws.Range("K12").Value = -19.23066314620145
$ws.Range("B13").Value = -19.23066314620145
$ws.Range("C13").Value = -19.23066314620145
$ws.Range("D13").Value = -19.23066314620145
$ws.Range("E13").Value = 2.52262420859753
$ws.Range("F13").Value = -19.23066314620145
$ws.Range("G13").Value = -19.23066314620145
$ws.Range("H13").Value = -19.23066314620145
$ws.Range("I13").Value = -19.23066314620145
$ws.Range("J13").Value = 1.75916718878794
$ws.Range("K13").Value = 1.801829935646107
$ws.Range("B14").Value = -19.23066314620145
$ws.Range("C14").Value = -19.23066314620145
$ws.Range("D14").Value = 1.564337255961574
$ws.Range("E14").Value = -19.23066314620145
$ws.Range("F14").Value = -19.23066314620145
$ws.Range("G14").Value = -19.23066314620145
$ws.Range("H14").Value = -19.23066314620145
$ws.Range("I14").Value = -19.23066314620145
$ws.Range("J14").Value = -19.23066314620145
$ws.Range("K14").Value = 1.976473343415951
$ws.Range("B15").Value = -19.23066314620145
$ws.Range("C15").Value = -19.23066314620145
$ws.Range("D15").Value = 1.742331359451139
$ws.Range("E15").Value = -19.23066314620145
$ws.Range("F15").Value = -19.23066314620145
$ws.Range("G15").Value = -19.23066314620145
$ws.Range("H15").Value = -19.23066314620145
$ws.Range("I15").Value = -19.23066314620145
$ws.Range("J15").Value = -19.23066314620145
$ws.Range("K15").Value = -19.23066314620145
$ws.Range("B16").Value = -19.23066314620145
$ws.Range("C16").Value = -19.23066314620145
$ws.Range("D16").Value = -19.23066314620145
$ws.Range("E16").Value = -19.23066314620145
$ws.Range("F16").Value = -19.23066314620145
$ws.Range("G16").Value = -19.23066314620145
$ws.Range("H16").Value = -19.23066314620145
$ws.Range("I16").Value = -19.23066314620145
$ws.Range("J16").Value = 1.905622726982808
$ws.Range("K16").Value = -19.23066314620145
$ws.Range("B17").Value = -19.23066314620145
$ws.Range("C17").Value = 2.122492834845521
$ws.Range("D17").Value = 1.803685145545166
$ws.Range("E17").Value = -19.23066314620145
$ws.Range("F17").Value = -19.23066314620145
$ws.Range("G17").Value = -19.23066314620145
$ws.Range("H17").Value = 2.031590591183894
$ws.Range("I17").Value = -19.23066314620145
$ws.Range("J17").Value = 2.427998567953601
$ws.Range("K17").Value = -19.23066314620145
$ws.Range("B18").Value = -19.23066314620145
$ws.Range("C18").Value = -19.23066314620145
$ws.Range("D18").Value = -19.23066314620145
$ws.Range("E18").Value = -19.23066314620145
$ws.Range("F18").Value = -19.23066314620145
$ws.Range("G18").Value = -19.23066314620145
$ws.Range("H18").Value = 1.985431993339576
$ws.Range("I18").Value = 4.321925867006125
$ws.Range("J18").Value = 2.41553384993355
$ws.Range("K18").Value = -19.23066314620145
$ws.Range("B19").Value = -19.23066314620145
$ws.Range("C19").Value = -19.23066314620145
$ws.Range("D19").Value = 2.045529997236425
$ws.Range("E19").Value = -19.23066314620145
$ws.Range("F19").Value = -19.23066314620145
$ws.Range("G19").Value = -19.23066314620145
$ws.Range("H19").Value = 1.611776976106998
$ws.Range("I19").Value = -19.23066314620145
$ws.Range("J19").Value = -19.23066314620145
$ws.Range("K19").Value = -19.23066314620145
$ws.Range("B20").Value = -19.23066314620145
$ws.Range("C20").Value = 1.075083569216863
$ws.Range("D20").Value = 1.517940022834853
$ws.Range("E20").Value = -19.23066314620145
$ws.Range("F20").Value = 3.232414396648644
$ws.Range("G20").Value = -19.23066314620145
$ws.Range("H20").Value = 1.497254816233637
$ws.Range("I20").Value = -19.23066314620145
$ws.Range("J20").Value = -19.23066314620145
$ws.Range("K20").Value = 2.114640843730599
$ws.Range("B21").Value = -19.23066314620145
$ws.Range("C21").Value = 1.309110026239585
$ws.Range("D21").Value = -19.23066314620145
$ws.Range("E21").Value = 1.710382310266914
$ws.Range("F21").Value = -19.23066314620145
$ws.Range("G21").Value = 2.533155852597223
$ws.Range("H21").Value = 1.59120501761585
$ws.Range("I21").Value = -19.23066314620145
$ws.Range("J21").Value = -19.23066314620145
$ws.Range("K21").Value = -19.23066314620145
